$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: in the opening "EXCELENTISSIMO ... COMARCA DE CONCHAS/SP." heading,
# underline just the word "COMARCA" (it keeps its existing bold formatting).
# "COMARCA" (upper case) occurs exactly once in the document, so a plain
# Find (no replace) safely locates the right run; we then flip the Font on
# the found Range, which splits the original run into three (pre / word /
# post) without disturbing anything else.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
[void]$rng1.Find.Execute("COMARCA", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng1.Find.Found) {
    $rng1.Font.Underline = 1
}

# ---------------------------------------------------------------------------
# Edit 2: "... no município de Pereiras, nesta comarca de Conchas, JOSÉ ..."
# becomes
#         "... no município de cidade, nesta comarca de comarca, JOSÉ ..."
# with the new second "comarca" underlined (just like "cidade" keeps the
# underline that "Pereiras" already had).
# ---------------------------------------------------------------------------

# 2a. "Pereiras" -> "cidade", preserving that run's underline formatting.
#     wdReplaceOne (the final argument, 1) is important: the document has a
#     second, unrelated "Pereiras" later on (Delegacia de Polícia de
#     Pereiras) that must stay untouched, and Find starting at the top of
#     $d.Content will hit the correct (first) occurrence.
$rng2 = $d.Content
[void]$rng2.Find.Execute("Pereiras", $true, $false, $false, $false, $false, $true, 1, $false, "cidade", 1)

# 2b. ", nesta comarca de Conchas, " -> ", nesta comarca de comarca, "
#     then underline just the newly introduced second "comarca".
$rng3 = $d.Content
$oldText = ", nesta comarca de Conchas, "
$newText = ", nesta comarca de comarca, "
[void]$rng3.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
if ($rng3.Find.Found) {
    $matchStart = $rng3.Start
    $wordOffset = $newText.LastIndexOf("comarca")
    $wordStart = $matchStart + $wordOffset
    $wordEnd = $wordStart + "comarca".Length
    $target = $d.Range($wordStart, $wordEnd)
    $target.Font.Underline = 1
}
